# Update the "取得日時" (retrieved timestamp) column on the "ランサーズ" sheet
# for rows 2-8 from 2026-01-09 12:40:16 to 2026-01-09 12:53:27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-01-09 12:53:27"
}
